# Update "想去人数" (F column, and in one case the min-price G column)
# figures for 展览 (sheet "展览") and 全部类型 (sheet "全部类型") worksheets,
# to reflect refreshed scrape numbers (output regenerated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (index 1) -------------------------------------------
$ws1 = $wb.Worksheets.Item("展览")

$ws1.Range("F2").Value  = 4726
$ws1.Range("F3").Value  = 1881
$ws1.Range("F4").Value  = 156
$ws1.Range("F6").Value  = 3182
$ws1.Range("F8").Value  = 599
$ws1.Range("F9").Value  = 287
$ws1.Range("F10").Value = 653
$ws1.Range("F11").Value = 553
$ws1.Range("F12").Value = 551
$ws1.Range("G12").Value = 50
$ws1.Range("F13").Value = 408
$ws1.Range("F14").Value = 144
$ws1.Range("F15").Value = 1798
$ws1.Range("F16").Value = 1390
$ws1.Range("F17").Value = 128
$ws1.Range("F18").Value = 1651
$ws1.Range("F20").Value = 131
$ws1.Range("F32").Value = 4021
$ws1.Range("F36").Value = 1579
$ws1.Range("F38").Value = 1903

# --- Sheet "全部类型" (index 4) -----------------------------------------
$ws4 = $wb.Worksheets.Item("全部类型")

$ws4.Range("F2").Value  = 4726
$ws4.Range("F3").Value  = 1881
$ws4.Range("F4").Value  = 156
$ws4.Range("F6").Value  = 3182
$ws4.Range("F8").Value  = 599
$ws4.Range("F9").Value  = 287
$ws4.Range("F10").Value = 653
$ws4.Range("F11").Value = 554
$ws4.Range("F12").Value = 551
$ws4.Range("G12").Value = 50
$ws4.Range("F14").Value = 408
$ws4.Range("F15").Value = 144
$ws4.Range("F16").Value = 1798
$ws4.Range("F17").Value = 1390
$ws4.Range("F18").Value = 128
$ws4.Range("F19").Value = 1651
$ws4.Range("F21").Value = 131
$ws4.Range("F33").Value = 4021
$ws4.Range("F39").Value = 1579
$ws4.Range("F41").Value = 1903
